$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "SV đạt" (E column) scores that changed
$ws.Range("E3").Value = 10
$ws.Range("E5").Value = 3
$ws.Range("E8").Value = 15

# Total row: turn the literal total into a real SUM formula
$ws.Range("E13").Formula = "=SUM(E3:E12)"

# Page setup: explicitly set to portrait orientation
$ws.PageSetup.Orientation = 1

# Move the active selection to E11 (single cell)
$ws.Range("E11").Select()
